$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for A2:D10 (A column = element number, B = Ni, C = Nj, D = t)
$data = @(
    @(0, 0, 6, 1),
    @(1, 2, 6, 1),
    @(2, 3, 8, 2),
    @(3, 4, 8, 2),
    @(4, 5, 7, 1),
    @(5, 1, 7, 1),
    @(6, 7, 9, 1),
    @(7, 6, 9, 1),
    @(8, 8, 9, 1)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $row++
}

# Update the selected cell to match the new selection in the diff
$ws.Range("H13").Select()
